# "Generate Report for Archive"
#
# 1. The localization status for the two source files has moved on from
#    "Ready for handoff" to "In Translation" - update every cell that
#    shows that status (the Overview sheet's per-language status columns,
#    plus the "Status" column on each per-language detail sheet).
# 2. The now-shorter status text means the "Status" columns can be
#    narrower - shrink them to fit.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find($oldStatus)
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        while ($true) {
            $found.Value = $newStatus
            $found = $used.FindNext($found)
            if ($found -eq $null -or $found.Address() -eq $firstAddress) {
                break
            }
        }
    }
}

# Narrow the "Status" columns now that the text is shorter.
# (ColumnWidth is expressed in characters and gets snapped by Excel to
# whole-pixel increments, so we dial in the input that lands on the
# nearest achievable width to the target.)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
